# Update the "Forecast Comparison" sheet with a new Week_Start_Date column
# and corrected week labels (strip leading zero from W01-W09 -> W1-W9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column B ("Week_Start_Date") - this shifts ASIN, MyForecast,
# Amazon Mean/P70/P80/P90 Forecast, Product Title, is_holiday_week one
# column to the right (B->C, C->D, D->E, E->F, F->G, G->H, H->I, I->J).
$ws.Columns("B:B").Insert()

# New header for the inserted column
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Make sure the new column stores the dates as plain text (not as Excel
# date serial numbers), matching the source data which used inline text.
$ws.Range("B2:B17").NumberFormat = "@"

$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
}

# Fix week labels in column A: drop the leading zero (W01 -> W1 .. W09 -> W9).
# W10 .. W16 remain unchanged.
for ($i = 1; $i -le 9; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "W$i"
}
